$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the dataset. It belongs right
# after the existing row 65 (chronologically among the "Femacal de La
# Calera" / "Bruselas (repollito)" records), so insert a fresh row at
# position 66 and push every row below it down by one (66->67, ..., 77->78).
$ws.Rows.Item(66).Insert()

# Populate the newly inserted row 66 with the new record's data.
$ws.Range("A66").Value = 3
$ws.Range("B66").Value = "Femacal de La Calera"
$ws.Range("C66").Value = "Coquimbo"
$ws.Range("D66").Value = 44798
$ws.Range("E66").Value = 5
$ws.Range("F66").Value = 100112035
$ws.Range("G66").Value = "Bruselas (repollito)"
$ws.Range("H66").Value = "Sin especificar"
$ws.Range("I66").Value = "Primera"
$ws.Range("J66").Value = 100
$ws.Range("K66").Value = 14000
$ws.Range("L66").Value = 15000
$ws.Range("M66").Value = 14450
$ws.Range("N66").Value = '$/malla 15 kilos'
$ws.Range("O66").Value = "Provincia de Quillota"
$ws.Range("P66").Value = 963
$ws.Range("Q66").Value = 15
$ws.Range("R66").Value = "Hortaliza"
